$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the A1:U5 block down to A8:U12 (values + existing formatting)
$ws.Range("A1:U5").Copy()
$ws.Range("A8").PasteSpecial(-4104)

# The first (header) row of the duplicated block (row 8) and the leading
# index column of the duplicated data rows (A9:A12) get a distinct look:
# same fill/border as before, but a lighter ("white, darker 15%") font.
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Font.ThemeColor = 2

$ws.Range("A1:U1").Copy()
$ws.Range("A8:U8").PasteSpecial(-4122)
$ws.Range("A8:U8").Font.ThemeColor = 2

$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)

# A lone, otherwise-empty cell further down/right picks up the sheet's
# highlight fill as well, extending the used range to Y15.
$ws.Range("A1").Copy()
$ws.Range("Y15").PasteSpecial(-4122)
$ws.Range("Y15").Font.Name = "Calibri"

$excel.CutCopyMode = 0

$ws.Range("W14").Select()
